$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as text so numeric-looking
# strings (e.g. "329.22") are not coerced into numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.063.30"
$ws.Range("E2").Value = "  -4.92%  "
$ws.Range("D3").Value = "1.828.83"
$ws.Range("E3").Value = "  -3.82%  "
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").Value = "329.22"
$ws.Range("E5").Value = "  -2.93%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").Value = "0.4641"
$ws.Range("E7").Value = "  -2.38%  "
$ws.Range("D8").Value = "0.3872"
$ws.Range("E8").Value = "  -3.43%  "
$ws.Range("D9").Value = "46.22"
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("D10").Value = "0.07878"
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("D11").Value = "0.9611"
$ws.Range("E11").Value = "  -3.13%  "
$ws.Range("D12").Value = "21.95"
$ws.Range("E12").Value = "  -5.58%  "
$ws.Range("D13").Value = "1.862.49"
$ws.Range("E13").Value = "  -2.28%  "
$ws.Range("D14").Value = "5.662"
$ws.Range("E14").Value = "  -4.81%  "
$ws.Range("D15").Value = "6.897"
$ws.Range("E15").Value = "  -2.95%  "
$ws.Range("D16").Value = "0.06861"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").Value = "0.9994"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "86.56"
$ws.Range("E18").Value = "  -2.88%  "
$ws.Range("D19").Value = "0.000009993"
$ws.Range("E19").Value = "  -2.04%  "
$ws.Range("D20").Value = "16.69"
$ws.Range("E20").Value = "  -3.81%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "28.094.51"
$ws.Range("E22").Value = "  -4.87%  "
$ws.Range("D23").Value = "5.322"
$ws.Range("E23").Value = "  -3.51%  "
$ws.Range("D24").Value = "11.02"
$ws.Range("E24").Value = "  -5.23%  "
$ws.Range("E25").Value = "  -2.86%  "
$ws.Range("D26").Value = "2.045.66"
$ws.Range("E26").Value = "  -4.07%  "
$ws.Range("D27").Value = "152.34"
$ws.Range("E27").Value = "  -3.13%  "
$ws.Range("D28").Value = "19.22"
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("D29").Value = "5.770"
$ws.Range("E29").Value = "  -11.43%  "
$ws.Range("D30").Value = "1.975"
$ws.Range("E30").Value = "  -3.95%  "
$ws.Range("D31").Value = "116.95"
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("D32").Value = "0.9395"
$ws.Range("E32").Value = "  -5.80%  "
$ws.Range("D33").Value = "0.09244"
$ws.Range("E33").Value = "  -3.09%  "
$ws.Range("D34").Value = "5.293"
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("D35").Value = "1.319"
$ws.Range("E35").Value = "  -5.02%  "
$ws.Range("D36").Value = "3.341"
$ws.Range("E36").Value = "  -5.39%  "
$ws.Range("E37").Value = "  -7.17%  "
$ws.Range("D38").Value = "0.02143"
$ws.Range("E38").Value = "  -4.68%  "
$ws.Range("D39").Value = "1.149"
$ws.Range("E39").Value = "  -4.30%  "
$ws.Range("D40").Value = "0.9996"
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("D41").Value = "7.648"
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("D42").Value = "0.5596"
$ws.Range("E42").Value = "  -4.01%  "
$ws.Range("D43").Value = "9.929"
$ws.Range("E43").Value = "  -5.96%  "
$ws.Range("D44").Value = "0.1768"
$ws.Range("E44").Value = "  -2.85%  "
$ws.Range("D45").Value = "1.204"
$ws.Range("E45").Value = "  -5.00%  "
$ws.Range("D46").Value = "2.233"
$ws.Range("E46").Value = "  -8.40%  "
$ws.Range("D47").Value = "11.60"
$ws.Range("E47").Value = "  -4.53%  "
$ws.Range("D48").Value = "0.5273"
$ws.Range("E48").Value = "  -4.06%  "
$ws.Range("D49").Value = "0.07003"
$ws.Range("E49").Value = "  -4.88%  "
$ws.Range("D50").Value = "1.829"
$ws.Range("E50").Value = "  -6.53%  "
$ws.Range("D51").Value = "111.49"
$ws.Range("E51").Value = "  -4.25%  "

# Remove the temporary number format so the cell keeps its original
# (unstyled) appearance, matching the source workbook.
$priceRange.ClearFormats()
